$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 corresponds to ano/ano_obj = 2025
# total_customers (C6): 415 -> 416
$ws.Range("C6").Value = 416

# new_customers (E6): 107 -> 108
$ws.Range("E6").Value = 108

# new_rate (G6): recalculated as new_customers / total_customers * 100
$ws.Range("G6").Value = 25.96153846153846

# returning_rate (H6): recalculated as returning_customers / total_customers * 100
$ws.Range("H6").Value = 74.03846153846155
